$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Feuil1")

$ws.Range("B2").Value = "17-77"
$ws.Range("B3").Value = "16-21"
$ws.Range("B4").Value = "25-39"
$ws.Range("B5").Value = "59-129"
$ws.Range("B6").Value = "133-135"
$ws.Range("B7").Value = "16-59"
$ws.Range("B8").Value = "89-102"
$ws.Range("B9").Value = "63-68"
$ws.Range("B10").Value = "70-76"
$ws.Range("B11").Value = "78-84"
$ws.Range("B12").Value = "105-119"
$ws.Range("B13").Value = "123-160"
$ws.Range("C13").Value = "envoiFormulaire"
$ws.Range("B14").Value = "1-6"
$ws.Range("B15").Value = "8-19"
$ws.Range("B16").Value = "23-41"

$ws.Range("B16").Select()

$wb.Save()
